$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("G2").Value = 1.030436666666667
$ws.Range("H2").Value = 3.09131
$ws.Range("O2").Value = 0.03059328965493693
$ws.Range("P2").Value = 0.03059328965493693
$ws.Range("Q2").Value = 0.1424891257455556
$ws.Range("R2").Value = 1.28240213171
$ws.Range("S2").Value = 0.03059328965493693
$ws.Range("T2").Value = 0.03059328965493693

# --- Row 3 updates ---
$ws.Range("G3").Value = 1.030436666666667
$ws.Range("H3").Value = 3.09131
$ws.Range("O3").Value = 0.1541543653555945
$ws.Range("P3").Value = 0.1541543653555945
$ws.Range("Q3").Value = 0.71797838667
$ws.Range("R3").Value = 6.461805480030001
$ws.Range("S3").Value = 0.1541543653555945
$ws.Range("T3").Value = 0.1541543653555945

# --- Row 4 updates ---
$ws.Range("G4").Value = 1.030436666666667
$ws.Range("H4").Value = 3.09131
$ws.Range("M4").Value = 3.682798
$ws.Range("N4").Value = 11.048394
$ws.Range("O4").Value = 0.8147861900435764
$ws.Range("P4").Value = 0.8147861900435764
$ws.Range("Q4").Value = 3.794890095126667
$ws.Range("R4").Value = 34.15401085614
$ws.Range("S4").Value = 0.8147861900435764
$ws.Range("T4").Value = 0.8147861900435764

# --- New row 5 ---
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Vip"
$ws.Range("C5").Value = "Vipr2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.030436666666667
$ws.Range("H5").Value = 3.09131
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.002107
$ws.Range("N5").Value = 0.006321
$ws.Range("O5").Value = 0.0004661549458921764
$ws.Range("P5").Value = 0.0004661549458921764
$ws.Range("Q5").Value = 0.002171130056666666
$ws.Range("R5").Value = 0.01954017051
$ws.Range("S5").Value = 0.0004661549458921764
$ws.Range("T5").Value = 0.0004661549458921764
